$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 398; this shifts existing rows 398-427 down to 399-428
$ws.Rows.Item(398).Insert()

# Populate the new row 398 with the new data record
$ws.Cells.Item(398, 1).Value = 10
$ws.Cells.Item(398, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(398, 3).Value = "La Araucanía"
$ws.Cells.Item(398, 4).Value = 45021
$ws.Cells.Item(398, 4).NumberFormat = $ws.Cells.Item(399, 4).NumberFormat
$ws.Cells.Item(398, 5).Value = 9
$ws.Cells.Item(398, 6).Value = 100112001
$ws.Cells.Item(398, 7).Value = "Berenjena"
$ws.Cells.Item(398, 8).Value = "Sin especificar"
$ws.Cells.Item(398, 9).Value = "Primera"
$ws.Cells.Item(398, 10).Value = 35
$ws.Cells.Item(398, 11).Value = 14000
$ws.Cells.Item(398, 12).Value = 14000
$ws.Cells.Item(398, 13).Value = 14000
$ws.Cells.Item(398, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(398, 15).Value = "Región del Maule"
$ws.Cells.Item(398, 16).Value = 350
$ws.Cells.Item(398, 17).Value = 40
$ws.Cells.Item(398, 18).Value = "Hortaliza"
